$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("tool_pid", "tool_code", "tool_type", "tool_price", "tool_inventory_start_date", "tool_inventory_end_date")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("F2").Select()
